# Applies the "change st_path name rule" edit described by the diff.
# Each call re-runs Find/Execute against the full document Content range
# so that every matching occurrence is replaced (Replace:=2 => wdReplaceAll).

$d = $word.ActiveDocument

function Replace-All($find, $replace, [bool]$matchCase, [bool]$matchWholeWord) {
    # NOTE: this engine does not honour PowerShell default parameter
    # values, so every caller below passes matchCase/matchWholeWord
    # explicitly.
    $d.Content.Find.Execute($find, $matchCase, $matchWholeWord, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# Simple global text replacements (old -> new), applied in an order that
# avoids any cascading between rules (see "08 Februari 2018" note below).
# MatchCase is always $true so the all-caps and mixed-case "BPS ..." runs
# stay distinct.

Replace-All "dgdb" "7778" $true $false
Replace-All "BPS PROVINSI SULAWESI TENGGARA" "BPS KABUPATEN BUTON" $true $false
Replace-All "RUBEN RAMBO" "AMALUDDIN HASAN" $true $false
Replace-All "196007251989031002" "196803091990031002" $true $false
Replace-All "Staf Seksi Statistik Pertanian" "KSK Lakudo" $true $false
Replace-All "yjdhfh" "dghs" $true $false
Replace-All "Wakatobi" "ddd" $true $false

# Date fields: "08 Februari 2018" must turn into "17 Februari 2018" BEFORE
# "06 Februari 2018" is turned into "08 Februari 2018", otherwise the
# freshly produced "08 Februari 2018" text would be caught by the first
# rule on a later pass.
Replace-All "08 Februari 2018" "17 Februari 2018" $true $false
Replace-All "07 Februari 2018" "16 Februari 2018" $true $false
Replace-All "06 Februari 2018" "08 Februari 2018" $true $false

Replace-All "Kendari" "dd" $true $false
Replace-All "Ir. H. ATQO MARDIYANTO, M.Si." "LA ODE MUSARAFA, SE" $true $false
Replace-All "196405081987021002" "196112311986031034" $true $false
Replace-All "fdhffhf" "fdgd" $true $false
Replace-All "DANI JAELANI, S.Si., MT" "SUDARWO" $true $false
Replace-All "BPS Provinsi Sulawesi Tenggara" "BPS Kabupaten Buton" $true $false
Replace-All "pesawat" "kapal" $true $false
Replace-All "196912101991121001" "198508292009011006" $true $false

# Budget / account code block (four separate paragraphs in the same cell).
Replace-All "054.01.06" "054.01.02" $true $false
Replace-All "2895" "2891" $true $false
Replace-All "027" "951" $true $false
Replace-All "062" "004" $true $false
Replace-All "524111" "524113" $true $false

# Single classification letter "C" -> "B" (unique run in the document).
Replace-All "C" "B" $true $true

Write-Host "edit applied"
